$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# NOTE: the order in which new text values are assigned below matches the
# order new shared-string entries should be appended in, so that the
# resulting xl/sharedStrings.xml table lines up with the target workbook.

# --- Line item table updates ---

# Row 14: single remaining line item
$ws.Range("E14").Value = "Airconditioner, Window type Brand: Kolin, Model: KAG-100HME4 Nominal capacity: 1.0 hp Cooling capacity: 9800kJ/h Refrigerant/Charge: R-410A / 460g Power Supply: 1PH / 230V / 60Hz Rated power: 865W Rated current: 3.8 A EER: 11.3"
$ws.Range("C14").Value = "unit"
$ws.Range("B14").Value = 1
$ws.Range("J14").Value = 45534

# --- Header / info block updates ---

# Row 7: Purchase Request (C7) and Department (I7)
$ws.Range("C7").Value = "Manila/Bacolod"
$ws.Range("I7").Value = "Electrical"

# Row 8: Date Prepared (C8) and Dept. Code (I8)
$ws.Range("C8").Value = 45498
$ws.Range("I8").Value = "EIC"

# Row 12: End-Use (C12)
$ws.Range("C12").Value = "SCADA ROOM"

# Row 11: Purpose (C11)
$ws.Range("C11").Value = "REPLACEMENT OF AIRCONDITIONING UNIT FOR SCADA ROOM"

# Row 9: Date Issued (C9) and Requestor (I9)
$ws.Range("C9").Value = 45498
$ws.Range("I9").Value = "GG"

# Row 10: Urgency No. (I10) cleared
$ws.Range("I10:K10").ClearContents()

# Rows 15-17: remove the other line items, unmerge and blank them out
$ws.Range("A15:K17").UnMerge()
$ws.Range("A15:K17").ClearContents()
$ws.Range("A15:K17").Style = "Normal"

# --- View changes ---
$ws.Range("J19").Select()
